# pulls r_2 values from excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("W3").Value = 0.9851912822914085
$ws.Range("Y3").Value = 0.01070685442323115

# Row 4 updates
$ws.Range("U4").Value = 186.7832931399875
$ws.Range("V4").Value = 0.01188074390501766
$ws.Range("W4").Value = 0.9850050806033371
